$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the "Increase/Decrease" column (W) from "Increase" to "Decrease"
# and the "Request Type" column (Y) from "Entry" to "Adjust" for rows 2-9.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 23).Value = "Decrease"   # Column W
    $ws.Cells.Item($r, 25).Value = "Adjust"     # Column Y
}

# Scroll the view so column P is the left-most visible column, then select X13.
$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("X13").Select()
